$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $n = $row - 1
    $ws.Range("C$row").Value = "Infocasas – Publicación $n"
    $ws.Range("D$row").Value = "Infocasas – Descripción $n"
}
